$wb = $excel.ActiveWorkbook

# Updated "想去人数" (want-to-go count) values for the 展览 and 全部类型 sheets.
$updates = @{
    2  = 7738
    3  = 7551
    4  = 108
    9  = 101
    11 = 217
    13 = 681
    14 = 1124
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
